$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.215.20'
$ws.Range('E2').Value = '  -0.22%  '
$ws.Range('D3').Value = '1.682.42'
$ws.Range('E3').Value = '  +0.25%  '
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '216.11'
$ws.Range('E5').Value = '  -0.67%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5254'
$ws.Range('E6').Value = '  -1.59%  '
$ws.Range('E7').Value = '  -0.07%  '
$ws.Range('E8').Value = '  +0.42%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06366'
$ws.Range('E9').Value = '  -1.56%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '21.45'
$ws.Range('E10').Value = '  -1.96%  '
$ws.Range('E11').Value = '  +1.20%  '
$ws.Range('D12').Value = '1.691.90'
$ws.Range('E12').Value = '  +0.78%  '
$ws.Range('E13').Value = '  -0.05%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.5746'
$ws.Range('E14').Value = '  -0.34%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.000008260'
$ws.Range('E15').Value = '  -2.38%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '66.09'
$ws.Range('E16').Value = '  +2.23%  '
$ws.Range('D17').Value = '26.246.48'
$ws.Range('E17').Value = '  -0.16%  '
$ws.Range('E18').Value = '  -0.05%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.873'
$ws.Range('E19').Value = '  -0.58%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '10.77'
$ws.Range('E20').Value = '  -0.76%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '189.37'
$ws.Range('E21').Value = '  -0.38%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.246'
$ws.Range('E22').Value = '  +0.75%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.006'
$ws.Range('E23').Value = '  -0.11%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '148.81'
$ws.Range('E24').Value = '  +2.21%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '7.756'
$ws.Range('E25').Value = '  -0.78%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.1259'
$ws.Range('E26').Value = '  -0.55%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '15.79'
$ws.Range('E27').Value = '  +0.36%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.06311'
$ws.Range('E28').Value = '  -2.47%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.379'
$ws.Range('E29').Value = '  -0.52%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.316'
$ws.Range('E30').Value = '  -0.23%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.568'
$ws.Range('E31').Value = '  -0.38%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.570'
$ws.Range('E32').Value = '  -0.49%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.681'
$ws.Range('E33').Value = '  +1.42%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.023'
$ws.Range('E34').Value = '  -0.69%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.6113'
$ws.Range('E35').Value = '  -1.15%  '
$ws.Range('E36').Value = '  +0.89%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.749'
$ws.Range('E37').Value = '  +1.17%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '6.171'
$ws.Range('E38').Value = '  -1.12%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01615'
$ws.Range('E39').Value = '  -0.45%  '
$ws.Range('D40').Value = '1.097.86'
$ws.Range('E40').Value = '  -1.21%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.8851'
$ws.Range('E41').Value = '  +1.48%  '
$ws.Range('E42').Value = '  -0.45%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '100.45'
$ws.Range('E43').Value = '  +0.14%  '
$ws.Range('D44').Value = '1.832.51'
$ws.Range('E44').Value = '  +0.31%  '
$ws.Range('E45').Value = '  +2.28%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '57.41'
$ws.Range('E46').Value = '  +0.60%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.004'
$ws.Range('E47').Value = '  +0.11%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '8.082'
$ws.Range('E48').Value = '  -0.83%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.05271'
$ws.Range('E49').Value = '  +0.17%  '
$ws.Range('E50').Value = '  -0.19%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '6.008'
$ws.Range('E51').Value = '  -1.20%  '
